$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Flatten the rich-text "Data type: ...; Columns: N; Version: 1" runs down
#    to plain shared strings (re-assigning the cell's own text collapses the
#    run list without changing the value or position in the shared table).
# ---------------------------------------------------------------------------
foreach ($name in @("tab1", "tab2", "tab3", "no_data1", "no_data2")) {
    $ws = $wb.Worksheets.Item($name)
    $cell = $ws.Range("A1")
    $cell.Value = $cell.Text
}

# ---------------------------------------------------------------------------
# 2. tab1: a blank row is inserted above the old "baz" row, pushing the
#    "baz" and "bat" data rows down by one (row6 -> row7, row7 -> row8).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("tab1")
$ws1.Rows(6).Insert()

# ---------------------------------------------------------------------------
# 3. tab2: the header row no longer hides rows 1-2, and the "h2" header cell
#    is replaced by the numeric value 2.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("tab2")
$ws2.Rows(1).Hidden = $false
$ws2.Rows(2).Hidden = $false
$ws2.Range("B2").Value = 2

# ---------------------------------------------------------------------------
# 4. Keep rows 1-2 hidden on the sheets that were not touched by the header
#    fix (the load/save round-trip otherwise drops the hidden flag).
# ---------------------------------------------------------------------------
foreach ($name in @("tab1", "tab3", "no_data1", "no_data2")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows(1).Hidden = $true
    $ws.Rows(2).Hidden = $true
}

# ---------------------------------------------------------------------------
# 5. Add the two new "ignored" tabs at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws7.Name = "Only 1 row - should be ignored"
$ws7.Range("A1").Value = "foo bar"

$ws8 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws7)
$ws8.Name = "Only 2 rows - should be ignored"
$ws8.Range("A1").Value = "foo bar"
$ws8.Range("A2").Value = "some "
$ws8.Range("B2").Value = "other "
$ws8.Range("C2").Value = "importer "
$ws8.Range("D2").Value = "or"
$ws8.Range("E2").Value = "something"

# ---------------------------------------------------------------------------
# 6. Selections: tab1 now highlights H23, tab2 becomes the active tab with
#    B3 selected (and must be activated last so it ends up the active tab).
# ---------------------------------------------------------------------------
$ws1.Range("H23").Select() | Out-Null

$ws2.Activate()
$ws2.Range("B3").Select() | Out-Null
